$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, pushing existing rows 8-93 down to 9-94.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data point.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44503
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112001
$ws.Range("G8").Value = "Berenjena"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8500
$ws.Range("M8").Value = 8250
$ws.Range("N8").Value = "$/caja 60 unidades"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 138
$ws.Range("Q8").Value = 60
$ws.Range("R8").Value = "Hortaliza"
